# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column D price cells to Text format (one cell at a time, since
# union/multi-area ranges only apply NumberFormat to the first area here)
# so numeric-looking strings (e.g. '46.00', '1.00', '33.867.61') are
# preserved exactly as text, matching the source inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "33.867.61"
$ws.Range("E2").Value = "  +9.62%  "
$ws.Range("D3").Value = "1.783.36"
$ws.Range("E3").Value = "  +6.37%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "224.91"
$ws.Range("E5").Value = "  +2.44%  "
$ws.Range("D6").Value = "0.558"
$ws.Range("E6").Value = "  +4.63%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "30.89"
$ws.Range("E8").Value = "  +6.17%  "
$ws.Range("D9").Value = "46.00"
$ws.Range("E9").Value = "  +4.28%  "
$ws.Range("D10").Value = "0.279"
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").Value = "0.0661"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "0.0923"
$ws.Range("E12").Value = "  +1.73%  "
$ws.Range("D13").Value = "2.040.57"
$ws.Range("E13").Value = "  +6.33%  "
$ws.Range("D14").Value = "1.778.41"
$ws.Range("D15").Value = "0.627"
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "33.839.58"
$ws.Range("E16").Value = "  +9.52%  "
$ws.Range("D17").Value = "10.01"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "4.20"
$ws.Range("E18").Value = "  +2.55%  "
$ws.Range("D19").Value = "68.65"
$ws.Range("E19").Value = "  +3.99%  "
$ws.Range("D20").Value = "250.93"
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").Value = "0.0₃0739"
$ws.Range("E21").Value = "  +2.83%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "10.29"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").Value = "4.22"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "2.16"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("D27").Value = "16.42"
$ws.Range("E27").Value = "  +3.90%  "
$ws.Range("D28").Value = "0.115"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").Value = "6.92"
$ws.Range("E29").Value = "  +3.82%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "3.81"
$ws.Range("E31").Value = "  +9.06%  "
$ws.Range("D32").Value = "0.0510"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").Value = "3.51"
$ws.Range("E34").Value = "  +5.86%  "
$ws.Range("D35").Value = "1.490.44"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +3.61%  "
$ws.Range("D38").Value = "0.626"
$ws.Range("E38").Value = "  +3.69%  "
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").Value = "82.94"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").Value = "2.35"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "0.886"
$ws.Range("E43").Value = "  +5.85%  "
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("D47").Value = "1.935.76"
$ws.Range("E47").Value = "  +6.72%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "5.73"
$ws.Range("E48").Value = "  +2.79%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "11.81"
$ws.Range("E50").Value = "  +11.94%  "
$ws.Range("D51").Value = "50.90"
$ws.Range("E51").Value = "  -0.60%  "
